$d = $word.ActiveDocument

# 1. Add ru-RU language to the "СТ РК 2864-2016" run
$rngLang = $d.Content
$null = $rngLang.Find.Execute("СТ РК 2864-2016", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngLang.LanguageID = "ru-RU"

# 2. Remove the _GoBack bookmark near the top of the document
$bmGoBack = $d.Bookmarks.Item("_GoBack")
$bmGoBack.Delete()

# 3. Replace the accepting_party date placeholder text, then split it into 3 runs: "{", "accepting_party_date", "}"
$rngA = $d.Content
$null = $rngA.Find.Execute('"{accepting_party_day}" {accepting_party_month} {accepting_party_year} г.', $false, $false, $false, $false, $false, $true, 1, $false, "{accepting_party_date}", 2)

$rngA2 = $d.Content
$null = $rngA2.Find.Execute("{accepting_party_date}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$aStart = $rngA2.Start
$aEnd = $rngA2.End
$aMiddle = $d.Range($aStart + 1, $aEnd - 1)
$aMiddle.Bold = 1
$aMiddle.Bold = 0

# 4. Replace the transferring_party date placeholder text, then split it into 3 runs: "{", "transferring_party_date", "}"
$rngT = $d.Content
$null = $rngT.Find.Execute('"{transferring_party_day}" {transferring_party_month} {transferring_party_year} г.', $false, $false, $false, $false, $false, $true, 1, $false, "{transferring_party_date}", 2)

$rngT2 = $d.Content
$null = $rngT2.Find.Execute("{transferring_party_date}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tStart = $rngT2.Start
$tEnd = $rngT2.End
$tMiddle = $d.Range($tStart + 1, $tEnd - 1)
$tMiddle.Bold = 1
$tMiddle.Bold = 0

# 5. Add a new _GoBack bookmark right after the second "_______________________" run
#    (the one in the "transferring party" / передающей cell), before the "(подпись)" run.
$rngLine = $d.Content
$null = $rngLine.Find.Execute("_______________________", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngLine2 = $d.Content
$rngLine2.Start = $rngLine.End
$null = $rngLine2.Find.Execute("_______________________", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertPoint = $d.Range($rngLine2.End, $rngLine2.End)
$d.Bookmarks.Add("_GoBack", $insertPoint)

Write-Output "done"
